# add properties for EInputEdit in QStudioSCADA and update ProjectManager's ui.
#
# The "变量显示文本框" (EInputEdit-related) entry row in the 画面编辑器 sheet
# is removed, and the active sheet / selection state is switched from the
# 事件功能 sheet back to 画面编辑器.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 画面编辑器
$ws2 = $wb.Worksheets.Item(2)   # 事件功能

# Remove the "变量显示文本框" / "未开始" row (row 12). A12 keeps its (blank)
# style slot, B12 is fully cleared (value + format) so it disappears.
$ws1.Range("A12").ClearContents()
$ws1.Range("B12").Clear()

# Switch the active sheet/selection back to 画面编辑器, with the cursor on
# B16, and drop 事件功能's tabSelected flag (its own selection stays put).
$null = $ws1.Activate()
$null = $ws1.Range("B16").Select()
